$d = $word.ActiveDocument

# 1) "BT=______2______" -> "BT=______1______"
$d.Content.Find.Execute("BT=______2______", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "BT=______1______", 2) | Out-Null

# 2) "AP=______126______" -> "AP=______211______"
$d.Content.Find.Execute("AP=______126______", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AP=______211______", 2) | Out-Null

# 3) "AR=_____126_______" -> "AR=_____211_______"
$d.Content.Find.Execute("AR=_____126_______", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AR=_____211_______", 2) | Out-Null
